# live_trading_results.xlsx - trade #67 closed update
# Commit: "Trade #67 closed at 2026-02-17 12:53:52 - unknown UNKNOWN +0.000%"
#
# 1) Summary sheet: refresh headline stats (capital, P&L, trade counts, win rate)
# 2) Strategy Status sheet: refresh the MarketMaking strategy row to match
# 3) All Trades / MarketMaking sheets: append the new closed-trade row (#67)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.16   # Current Capital
$summary.Range("B4").Value = 0.15      # Total P&L $
$summary.Range("B6").Value = 67        # Total Trades
$summary.Range("B7").Value = 30        # Winning Trades
$summary.Range("B9").Value = 44.78     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status (MarketMaking is row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.16     # Capital
$status.Range("D4").Value = 67         # Trades
$status.Range("E4").Value = 0.15       # P&L $
$status.Range("F4").Value = 0.16       # P&L %
$status.Range("G4").Value = 44.78      # Win Rate %

# ---------------------------------------------------------------------------
# 3) Append trade #67 to "All Trades" and "MarketMaking" sheets (new row 68)
# ---------------------------------------------------------------------------
$newRow = @{
    A = 67
    B = "2026-02-17"
    C = "12:53:45"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.91
    G = 0.92
    H = "CLOSED"
    I = 1.0989
    J = 0.01
    K = 100.16
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A68").Value = $newRow.A
    # Date-looking text needs a text prefix so Excel doesn't coerce it to a
    # date serial number (matches how a literal "2026-02-17" behaves when
    # typed into a General-formatted cell).
    $ws.Range("B68").Value = "'" + $newRow.B
    $ws.Range("C68").Value = $newRow.C
    $ws.Range("D68").Value = $newRow.D
    $ws.Range("E68").Value = $newRow.E
    $ws.Range("F68").Value = $newRow.F
    $ws.Range("G68").Value = $newRow.G
    $ws.Range("H68").Value = $newRow.H
    $ws.Range("I68").Value = $newRow.I
    $ws.Range("J68").Value = $newRow.J
    $ws.Range("K68").Value = $newRow.K
    $ws.Range("L68").Value = $newRow.L
    $ws.Range("M68").Value = $newRow.M
    $ws.Range("N68").Value = $newRow.N
    $ws.Range("O68").Value = $newRow.O
    $ws.Range("P68").Value = $newRow.P
    $ws.Range("Q68").Value = $newRow.Q
}
